$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 271 (shifts existing rows 271.. down by one)
$ws.Rows.Item(271).Insert()

# Populate the newly inserted row 271 with the new data record
$ws.Cells.Item(271, 1).Value = 3
$ws.Cells.Item(271, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(271, 3).Value = "Coquimbo"
$ws.Cells.Item(271, 4).Value = 44755
$ws.Cells.Item(271, 5).Value = 5
$ws.Cells.Item(271, 6).Value = 100114013
$ws.Cells.Item(271, 7).Value = "Zanahoria"
$ws.Cells.Item(271, 8).Value = "Sin especificar"
$ws.Cells.Item(271, 9).Value = "Primera"
$ws.Cells.Item(271, 10).Value = 210
$ws.Cells.Item(271, 11).Value = 11000
$ws.Cells.Item(271, 12).Value = 11000
$ws.Cells.Item(271, 13).Value = 11000
$ws.Cells.Item(271, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(271, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(271, 16).Value = 550
$ws.Cells.Item(271, 17).Value = 20
$ws.Cells.Item(271, 18).Value = "Hortaliza"

# Apply the date style (numFmt) used by the rest of column D to the new D271 cell
$ws.Cells.Item(271, 4).NumberFormat = $ws.Cells.Item(272, 4).NumberFormat
